# Apply cryptocurrency price/volume updates per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.401.10"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "1.825.47"
$ws.Range("E3").Value = "  -2.98%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$ws.Range("E4").Value = "  -0.97%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "330.58"

$ws.Range("E6").Value = "  -0.76%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4586"
$ws.Range("E7").Value = "  -2.25%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3814"
$ws.Range("E8").Value = "  -3.83%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.30"
$ws.Range("E9").Value = "  +1.43%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07896"
$ws.Range("E10").Value = "  -1.83%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.9647"
$ws.Range("E11").Value = "  -5.00%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "21.05"
$ws.Range("E12").Value = "  -4.56%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.876"
$ws.Range("E13").Value = "  -2.43%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.825.27"
$ws.Range("E14").Value = "  -2.92%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.093"
$ws.Range("E15").Value = "  -2.82%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.003"
$ws.Range("E16").Value = "  -0.86%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "89.73"
$ws.Range("E17").Value = "  +0.49%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06603"
$ws.Range("E18").Value = "  -1.87%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.00001023"
$ws.Range("E19").Value = "  -2.48%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.19"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("D22").Value = "27.395.98"
$ws.Range("E22").Value = "  -1.76%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.305"
$ws.Range("E23").Value = "  -3.70%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.85"
$ws.Range("E24").Value = "  -1.66%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.277"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("D26").Value = "2.048.45"
$ws.Range("E26").Value = "  -2.53%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "156.27"
$ws.Range("E27").Value = "  -2.05%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.41"
$ws.Range("E28").Value = "  -2.32%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.069"
$ws.Range("E29").Value = "  -4.24%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.276"
$ws.Range("E30").Value = "  -4.13%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "118.14"
$ws.Range("E31").Value = "  -3.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09322"
$ws.Range("E32").Value = "  -2.02%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9372"
$ws.Range("E33").Value = "  -5.21%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.572"
$ws.Range("E34").Value = "  -1.79%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.228"
$ws.Range("E35").Value = "  -2.47%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.322"
$ws.Range("E36").Value = "  -2.57%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.05927"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02181"
$ws.Range("E38").Value = "  -3.11%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.128"
$ws.Range("E39").Value = "  -2.29%  "

$ws.Range("E40").Value = "  -0.71%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.143"
$ws.Range("E41").Value = "  -4.67%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.5767"
$ws.Range("E42").Value = "  -4.18%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1826"
$ws.Range("E43").Value = "  -3.96%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "9.970"
$ws.Range("E44").Value = "  -4.23%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.264"
$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5422"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "11.83"
$ws.Range("E47").Value = "  -3.31%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.873"
$ws.Range("E48").Value = "  -3.78%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06565"
$ws.Range("E49").Value = "  -3.29%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "109.71"
$ws.Range("E50").Value = "  -2.69%  "

$ws.Range("E51").Value = "  -34.07%  "
